# Classe: OutlookMail - concluída
# Edits slide 9 ("Classe: OutlookMail") of the deck:
#   - Fix typo "SMPT Server" -> "SMP Server"
#   - Add a new bullet paragraph "SMTP Port" right after it
#   - Remove the "Retorna Log" bullet paragraph from the "Enviar e-mail" box
#   - Reposition the connector line / the two text boxes to the values
#     PowerPoint settles on after the text edits (autofit growth/shrink).

# Shape.Top/.Left are exposed as single-precision points, so a naive
# emu/12700.0 conversion can truncate to one EMU below the intended value
# once it is rounded back to EMU on save. Nudging by half an EMU before
# the float32 round-trip lands exactly on the target EMU value.
function EmuToPoints($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

$connector   = $s.Shapes.Item(3)   # "Conector reto 12"
$serverBox   = $s.Shapes.Item(5)   # "CaixaDeTexto 5" (Assunto/Origem/.../SMPT Server)
$sendBox     = $s.Shapes.Item(6)   # "CaixaDeTexto 6" (Enviar e-mail / Retorna Log)

# --- 1. "SMPT Server" -> "SMP Server" ------------------------------------
$tr = $serverBox.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf("SMPT Server")
$old = $tr.Characters($idx + 1, 11)
$old.Text = "SMP Server"

# --- 2. Append a new "SMTP Port" bullet paragraph ------------------------
$nl = [char]13
$null = $tr.InsertAfter($nl + "SMTP Port")

# --- 3. Remove the "Retorna Log" paragraph from the send box -------------
$tr2 = $sendBox.TextFrame.TextRange
$lastPara = $tr2.Paragraphs(2, 1)
$lastPara.Delete()
$lastParaEmpty = $tr2.Paragraphs(2, 1)
$lastParaEmpty.Delete()

# --- 4. Re-flow the vertical layout --------------------------------------
# (the text boxes use spAutoFit, so Height updates on its own; Top needs to
# be set explicitly to match the post-edit layout)
$connector.Top = EmuToPoints 5305645
$serverBox.Top = EmuToPoints 2837534
$sendBox.Top = EmuToPoints 5594128
